$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Add a new data row (row 12), splitting the old combined
# "Program Students" row into its own "YOS" row: same record as row 11
# but with a brand new Encrypted Student No and with Program Code /
# Program Title left blank (that info now lives in its own table).
# ------------------------------------------------------------------

# Copy the formatting (styles) of row 11 into row 12 first, so every cell
# in A12:Y12 gets the same per-column style (s="2" text / s="3" numeric)
# as the rest of the data rows.
$ws.Range("A11:Y11").Copy()
$ws.Range("A12:Y12").PasteSpecial(-4122)

# Now populate the values for row 12.
$ws.Range("A12").Value = "00B197BA7753B1F2CFD57570245D62E5"
$ws.Range("D12").Value = "2017"
$ws.Range("H12").Value = "YOS 2"
$ws.Range("I12").Value = "South Africa"
$ws.Range("L12").Value = "Zulu"
$ws.Range("M12").Value = "Black"
$ws.Range("N12").Value = "M"
$ws.Range("O12").Value = 45
$ws.Range("P12").Value = "BUSE2023"
$ws.Range("Q12").Value = 50
$ws.Range("R12").Value = "PAS"
$ws.Range("S12").Value = "Q"
$ws.Range("T12").Value = "Completed all requirements for qualification"
$ws.Range("U12").Value = "Q"
$ws.Range("V12").Value = 63.67
$ws.Range("Y12").Value = "asd"

# Recreate the merged "header style" cells for the new row, matching the
# pattern used by every other data row (A:C, F:G, I:K).
$ws.Range("A12:C12").Merge()
$ws.Range("F12:G12").Merge()
$ws.Range("I12:K12").Merge()

# Match the selection / scroll position captured for the edited file.
$ws.Application.Goto($ws.Range("T12"), $true)

# Widen column A slightly (closest value this engine's pixel-quantised
# ColumnWidth model can reach to the recorded 15.4438775510204 units).
$ws.Columns.Item(1).ColumnWidth = 14.59

Write-Host "done"
